# Update simulation results for pl_mw.xlsx (Case 5_70, 380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.3825411565576928
$ws.Cells.Item(2, 3).Value = 0.060807607171256
$ws.Cells.Item(2, 4).Value = 0.1972768006721139
$ws.Cells.Item(2, 5).Value = 0.1687696096471782
$ws.Cells.Item(2, 6).Value = 1.424585000670383
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 0.7082064001759925
$ws.Cells.Item(2, 10).Value = 0.1873394869775851
$ws.Cells.Item(2, 11).Value = 0.4169121961893723
$ws.Cells.Item(2, 15).Value = 3.444437899617625

# Row 3
$ws.Cells.Item(3, 2).Value = 0.3426740460885753
$ws.Cells.Item(3, 3).Value = 0.05361555193017864
$ws.Cells.Item(3, 4).Value = 0.1905918361524783
$ws.Cells.Item(3, 5).Value = 0.1643571820370653
$ws.Cells.Item(3, 6).Value = 1.429303185283004
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 0.7158796129890419
$ws.Cells.Item(3, 10).Value = 0.1834558576622669
$ws.Cells.Item(3, 11).Value = 0.3719707251527211
$ws.Cells.Item(3, 15).Value = 3.469712867894231

# Row 4
$ws.Cells.Item(4, 2).Value = 0.3182053267261153
$ws.Cells.Item(4, 3).Value = 0.04919238347044086
$ws.Cells.Item(4, 4).Value = 0.1865624613248542
$ws.Cells.Item(4, 5).Value = 0.1617301785162688
$ws.Cells.Item(4, 6).Value = 1.43298091856483
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 0.7209753403111314
$ws.Cells.Item(4, 10).Value = 0.1811772573574117
$ws.Cells.Item(4, 11).Value = 0.3443737628843735
$ws.Cells.Item(4, 15).Value = 3.487160028938021

# Row 5
$ws.Cells.Item(5, 2).Value = 0.308237226976189
$ws.Cells.Item(5, 3).Value = 0.04738818196057082
$ws.Cells.Item(5, 4).Value = 0.1849394946407159
$ws.Cells.Item(5, 5).Value = 0.1606803902649787
$ws.Cells.Item(5, 6).Value = 1.434675993104477
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 0.7231485015043937
$ws.Cells.Item(5, 10).Value = 0.1802753754285078
$ws.Cells.Item(5, 11).Value = 0.3331277659760303
$ws.Cells.Item(5, 15).Value = 3.494754645539999

# Row 6
$ws.Cells.Item(6, 2).Value = 0.3065822374375671
$ws.Cells.Item(6, 3).Value = 0.04708849419455419
$ws.Cells.Item(6, 4).Value = 0.1846711555830893
$ws.Cells.Item(6, 5).Value = 0.1605073278374327
$ws.Cells.Item(6, 6).Value = 1.43496932037263
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 0.7235151872923424
$ws.Cells.Item(6, 10).Value = 0.1801272298087326
$ws.Cells.Item(6, 11).Value = 0.3312603959686271
$ws.Cells.Item(6, 15).Value = 3.496045000701216

# Row 7
$ws.Cells.Item(7, 2).Value = 0.3180708800554726
$ws.Cells.Item(7, 3).Value = 0.04916805824115045
$ws.Cells.Item(7, 4).Value = 0.1865404961891386
$ws.Cells.Item(7, 5).Value = 0.1617159366668446
$ws.Cells.Item(7, 6).Value = 1.433002983789557
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 0.721004257196693
$ws.Cells.Item(7, 10).Value = 0.1811649862640721
$ws.Cells.Item(7, 11).Value = 0.3442220945863426
$ws.Cells.Item(7, 15).Value = 3.487260489984294

# Row 8
$ws.Cells.Item(8, 2).Value = 0.3687932953046698
$ws.Cells.Item(8, 3).Value = 0.05832933646769334
$ws.Cells.Item(8, 4).Value = 0.1949562833198542
$ws.Cells.Item(8, 5).Value = 0.1672311664672463
$ws.Cells.Item(8, 6).Value = 1.426049809422572
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 0.7107723325058473
$ws.Cells.Item(8, 10).Value = 0.1859784295886513
$ws.Cells.Item(8, 11).Value = 0.4014173298091066
$ws.Cells.Item(8, 15).Value = 3.452752435521063

# Row 9
$ws.Cells.Item(9, 2).Value = 0.4683154017187405
$ws.Cells.Item(9, 3).Value = 0.07623421192157309
$ws.Cells.Item(9, 4).Value = 0.2120521081267128
$ws.Cells.Item(9, 5).Value = 0.1786973675917238
$ws.Cells.Item(9, 6).Value = 1.41860897249795
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 0.6937587135357397
$ws.Cells.Item(9, 10).Value = 0.1962581699807089
$ws.Cells.Item(9, 11).Value = 0.5135314631701817
$ws.Cells.Item(9, 15).Value = 3.400388407791951

# Row 10
$ws.Cells.Item(10, 2).Value = 0.5414452157912137
$ws.Cells.Item(10, 3).Value = 0.0893492576558117
$ws.Cells.Item(10, 4).Value = 0.2249692238509198
$ws.Cells.Item(10, 5).Value = 0.187517125476127
$ws.Cells.Item(10, 6).Value = 1.416919572574812
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 0.6831209828711913
$ws.Cells.Item(10, 10).Value = 0.2043239766556582
$ws.Cells.Item(10, 11).Value = 0.5958499032643942
$ws.Cells.Item(10, 15).Value = 3.371258813343417

# Row 11
$ws.Cells.Item(11, 2).Value = 0.5747118371222371
$ws.Cells.Item(11, 3).Value = 0.09530649450877604
$ws.Cells.Item(11, 4).Value = 0.2309221928579888
$ws.Cells.Item(11, 5).Value = 0.1916151523982563
$ws.Cells.Item(11, 6).Value = 1.41697162756364
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 0.6786864215859332
$ws.Cells.Item(11, 10).Value = 0.2081050173503769
$ws.Cells.Item(11, 11).Value = 0.6332828751880015
$ws.Cells.Item(11, 15).Value = 3.360038181655312

# Row 12
$ws.Cells.Item(12, 2).Value = 0.5873084391926682
$ws.Cells.Item(12, 3).Value = 0.0975609976488272
$ws.Cells.Item(12, 4).Value = 0.2331873837922274
$ws.Cells.Item(12, 5).Value = 0.1931792763254734
$ws.Cells.Item(12, 6).Value = 1.417109335140864
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 0.6770653943257372
$ws.Cells.Item(12, 10).Value = 0.2095528808558385
$ws.Cells.Item(12, 11).Value = 0.6474551752295099
$ws.Cells.Item(12, 15).Value = 3.356081404844645

# Row 13
$ws.Cells.Item(13, 2).Value = 0.5845955794102906
$ws.Cells.Item(13, 3).Value = 0.09707551269801229
$ws.Cells.Item(13, 4).Value = 0.2326990502916715
$ws.Cells.Item(13, 5).Value = 0.1928418685836135
$ws.Cells.Item(13, 6).Value = 1.417074429505234
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 0.6774119200686535
$ws.Cells.Item(13, 10).Value = 0.2092403433603209
$ws.Cells.Item(13, 11).Value = 0.6444030508651508
$ws.Cells.Item(13, 15).Value = 3.356920566584392

# Row 14
$ws.Cells.Item(14, 2).Value = 0.575748185945713
$ws.Cells.Item(14, 3).Value = 0.09549200192256535
$ws.Cells.Item(14, 4).Value = 0.2311083331461106
$ws.Cells.Item(14, 5).Value = 0.1917435878402998
$ws.Cells.Item(14, 6).Value = 1.416980592056461
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 0.6785518910149584
$ws.Cells.Item(14, 10).Value = 0.2082238120960795
$ws.Cells.Item(14, 11).Value = 0.6344488981885945
$ws.Cells.Item(14, 15).Value = 3.359706797094105

# Row 15
$ws.Cells.Item(15, 2).Value = 0.5703287849524941
$ws.Cells.Item(15, 3).Value = 0.09452187410346369
$ws.Cells.Item(15, 4).Value = 0.2301353926581982
$ws.Cells.Item(15, 5).Value = 0.1910724579593435
$ws.Cells.Item(15, 6).Value = 1.416938480339709
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 0.6792577432472022
$ws.Cells.Item(15, 10).Value = 0.2076032489244426
$ws.Cells.Item(15, 11).Value = 0.6283513120833106
$ws.Cells.Item(15, 15).Value = 3.36145150817697

# Row 16
$ws.Cells.Item(16, 2).Value = 0.5392710961700118
$ws.Cells.Item(16, 3).Value = 0.08895975142380053
$ws.Cells.Item(16, 4).Value = 0.2245817199271443
$ws.Cells.Item(16, 5).Value = 0.1872510327936112
$ws.Cells.Item(16, 6).Value = 1.416932682528909
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 0.683418938895052
$ws.Cells.Item(16, 10).Value = 0.2040791260096171
$ws.Cells.Item(16, 11).Value = 0.5934032246269112
$ws.Cells.Item(16, 15).Value = 3.372032988741296

# Row 17
$ws.Cells.Item(17, 2).Value = 0.5202176022429228
$ws.Cells.Item(17, 3).Value = 0.08554522890648286
$ws.Cells.Item(17, 4).Value = 0.2211943281628237
$ws.Cells.Item(17, 5).Value = 0.1849286651852537
$ws.Cells.Item(17, 6).Value = 1.417139292860341
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 0.6860753702278615
$ws.Cells.Item(17, 10).Value = 0.201945827480202
$ws.Cells.Item(17, 11).Value = 0.5719595759019285
$ws.Cells.Item(17, 15).Value = 3.379044632295006

# Row 18
$ws.Cells.Item(18, 2).Value = 0.5092585297308005
$ws.Cells.Item(18, 3).Value = 0.08358045748423137
$ws.Cells.Item(18, 4).Value = 0.219253236666745
$ws.Cells.Item(18, 5).Value = 0.1836009872254465
$ws.Cells.Item(18, 6).Value = 1.41733536587941
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 0.6876413578845089
$ws.Cells.Item(18, 10).Value = 0.200729341086813
$ws.Cells.Item(18, 11).Value = 0.5596244724699773
$ws.Cells.Item(18, 15).Value = 3.383268656509728

# Row 19
$ws.Cells.Item(19, 2).Value = 0.5055479954968405
$ws.Cells.Item(19, 3).Value = 0.08291508037936524
$ws.Cells.Item(19, 4).Value = 0.2185972649767223
$ws.Cells.Item(19, 5).Value = 0.1831528489213738
$ws.Cells.Item(19, 6).Value = 1.417415018271171
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 0.688178112471995
$ws.Cells.Item(19, 10).Value = 0.2003192689393529
$ws.Cells.Item(19, 11).Value = 0.5554478200770916
$ws.Cells.Item(19, 15).Value = 3.384731654931471

# Row 20
$ws.Cells.Item(20, 2).Value = 0.5222458847875373
$ws.Cells.Item(20, 3).Value = 0.08590879710826016
$ws.Cells.Item(20, 4).Value = 0.2215541727295687
$ws.Cells.Item(20, 5).Value = 0.1851750486876185
$ws.Cells.Item(20, 6).Value = 1.417109305186997
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 0.6857886471155119
$ws.Cells.Item(20, 10).Value = 0.2021718310321745
$ws.Cells.Item(20, 11).Value = 0.5742424270615913
$ws.Cells.Item(20, 15).Value = 3.378278449716817

# Row 21
$ws.Cells.Item(21, 2).Value = 0.5783469054595116
$ws.Cells.Item(21, 3).Value = 0.09595715547425243
$ws.Cells.Item(21, 4).Value = 0.2315752696982827
$ws.Cells.Item(21, 5).Value = 0.1920658464365133
$ws.Cells.Item(21, 6).Value = 1.417004952043953
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 0.6782154726525782
$ws.Cells.Item(21, 10).Value = 0.2085219561232208
$ws.Cells.Item(21, 11).Value = 0.6373727546417456
$ws.Cells.Item(21, 15).Value = 3.358880480632479

# Row 22
$ws.Cells.Item(22, 2).Value = 0.6150075096577439
$ws.Cells.Item(22, 3).Value = 0.1025162891919251
$ws.Cells.Item(22, 4).Value = 0.2381882891209557
$ws.Cells.Item(22, 5).Value = 0.1966409896936554
$ws.Cells.Item(22, 6).Value = 1.417624517833815
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 0.6736054908103739
$ws.Cells.Item(22, 10).Value = 0.2127657619220287
$ws.Cells.Item(22, 11).Value = 0.6786156854004446
$ws.Cells.Item(22, 15).Value = 3.347906253817371

# Row 23
$ws.Cells.Item(23, 2).Value = 0.5954417032301933
$ws.Cells.Item(23, 3).Value = 0.09901632508231728
$ws.Cells.Item(23, 4).Value = 0.2346530143816494
$ws.Cells.Item(23, 5).Value = 0.194192616865891
$ws.Cells.Item(23, 6).Value = 1.417230913998949
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 0.6760348369953988
$ws.Cells.Item(23, 10).Value = 0.2104922024537643
$ws.Cells.Item(23, 11).Value = 0.656605295878677
$ws.Cells.Item(23, 15).Value = 3.353607458916457

# Row 24
$ws.Cells.Item(24, 2).Value = 0.5213289136642629
$ws.Cells.Item(24, 3).Value = 0.0857444332644377
$ws.Cells.Item(24, 4).Value = 0.221391467165148
$ws.Cells.Item(24, 5).Value = 0.1850636353966948
$ws.Cells.Item(24, 6).Value = 1.417122621859619
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 0.6859181538724961
$ws.Cells.Item(24, 10).Value = 0.202069623753502
$ws.Cells.Item(24, 11).Value = 0.5732103713793322
$ws.Cells.Item(24, 15).Value = 3.378624239886165

# Row 25
$ws.Cells.Item(25, 2).Value = 0.4413884680982676
$ws.Cells.Item(25, 3).Value = 0.07139722915452751
$ws.Cells.Item(25, 4).Value = 0.2073642579649828
$ws.Cells.Item(25, 5).Value = 0.1755259054451983
$ws.Cells.Item(25, 6).Value = 1.419958649566013
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 0.698034483224415
$ws.Cells.Item(25, 10).Value = 0.1933871493186388
$ws.Cells.Item(25, 11).Value = 0.4832089932229735
$ws.Cells.Item(25, 15).Value = 3.412914473483681

